$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top, pushing existing data down
$ws.Rows.Item(1).Insert()

# Fill in the header row
$headers = @("CNE", "FirstName", "LastName", "DateofBirth", "ClasseName", "Phone", "Email")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Match the style used elsewhere in the sheet (no fill/border) for header cells
$ws.Range("A1:G1").Style = "Normal"

# Update selection to mirror the saved view state
$ws.Range("G12").Select()
